# The old importer emitted the "<id>...</id>" marker as three separately
# formatted runs per occurrence:
#   run 1: "<id>"        (Courier New, brownish color, sz 18)
#   run 2: "<the id>"    (plain, black)
#   run 3: "</id>"       (Courier New, brownish color, sz 18)
# Newly (re-)downloaded markup instead represents the whole tag as a single
# run/string, e.g. "<id>p103r_1</id>". Collapse each such triple back into
# one run, keeping the character formatting of the opening "<id>" run
# (matches what Word does when you edit/retype across adjoining runs).
#
# We don't hardcode positions/ids - Find locates every "<id>" ... "</id>"
# pair in the document so this applies to however many occurrences exist.

$d = $word.ActiveDocument
$docEnd = $d.Content.End

$searchFrom = 0
$mergeCount = 0

while ($true) {
    $openRng = $d.Range($searchFrom, $docEnd)
    $openFound = $openRng.Find.Execute("<id>", $true, $false, $false, $false, $false, `
                                        $true, 1, $false, "", 0)
    if (-not $openFound) { break }

    $openStart = $openRng.Start
    $openEnd = $openRng.End

    $closeRng = $d.Range($openEnd, $docEnd)
    $closeFound = $closeRng.Find.Execute("</id>", $true, $false, $false, $false, $false, `
                                          $true, 1, $false, "", 0)
    if (-not $closeFound) { break }

    $closeStart = $closeRng.Start
    $closeEnd = $closeRng.End

    # Text sitting between the two tags is the id value itself.
    $midRng = $d.Range($openEnd, $closeStart)
    $idValue = $midRng.Text

    $replacement = $idValue + "</id>"

    # Remove the old "<id value>" + "</id>" runs entirely ...
    $spanRng = $d.Range($openEnd, $closeEnd)
    $spanRng.Delete()

    # ... and retype the same text right after the "<id>" run, so it
    # merges into that run and inherits its formatting.
    $openRunRng = $d.Range($openStart, $openEnd)
    $openRunRng.InsertAfter($replacement)

    $mergeCount = $mergeCount + 1

    $searchFrom = $openStart + 4 + $replacement.Length
    $docEnd = $d.Content.End
}

Write-Output "Merged $mergeCount <id>...</id> run group(s)."
